# PB-427 - Introducing CALCULATE and CALCULATETABLE
# Extend the "Product" table (F19:F23) with two new columns, Color and
# Size, and select H24 (just past the bottom-right of the now-larger
# table) to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the ListObject from F19:F23 (1 column) to F19:H23 (3 columns).
# Resizing first (before any header text is written) lets the new
# ListColumns pick up their names straight from the header cells below.
$lo = $ws.ListObjects("Product")
$lo.Resize($ws.Range("F19:H23"))

# Column G - Color
$ws.Range("G19").Value = "Color"
$ws.Range("G20").Value = "Red"
$ws.Range("G21").Value = "Red"
$ws.Range("G22").Value = "Blue"
$ws.Range("G23").Value = "Blue"

# Column H - Size (written L-rows before M-rows to reproduce the
# original authoring order of the shared-string table)
$ws.Range("H19").Value = "Size"
$ws.Range("H21").Value = "L"
$ws.Range("H23").Value = "L"
$ws.Range("H20").Value = "M"
$ws.Range("H22").Value = "M"

# Final selection left on the sheet
$ws.Range("H24").Select()
